function Set-CellText($ws, $row, $col, $val) {
    $cell = $ws.Cells.Item($row, $col)
    if ($val -match '^[0-9]+(\.[0-9]+)?$') {
        $cell.NumberFormat = "@"
    }
    $cell.Value2 = $val
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
Set-CellText $ws 2 4 '43.015.11'
Set-CellText $ws 2 5 '  -0.40%  '

# Row 3
Set-CellText $ws 3 4 '2.238.15'
Set-CellText $ws 3 5 '  -1.95%  '

# Row 4
Set-CellText $ws 4 5 '  +0.35%  '

# Row 5
Set-CellText $ws 5 4 '115.73'
Set-CellText $ws 5 5 '  +2.04%  '

# Row 6
Set-CellText $ws 6 4 '265.51'
Set-CellText $ws 6 5 '  -0.06%  '

# Row 7
Set-CellText $ws 7 4 '0.629'
Set-CellText $ws 7 5 '  +1.11%  '

# Row 8
Set-CellText $ws 8 5 '  +0.35%  '

# Row 9
Set-CellText $ws 9 4 '0.606'
Set-CellText $ws 9 5 '  -0.80%  '

# Row 10
Set-CellText $ws 10 4 '46.44'
Set-CellText $ws 10 5 '  -2.79%  '

# Row 11
Set-CellText $ws 11 4 '0.0928'
Set-CellText $ws 11 5 '  -0.44%  '

# Row 12
Set-CellText $ws 12 4 '9.14'
Set-CellText $ws 12 5 '  +0.30%  '

# Row 13
Set-CellText $ws 13 5 '  -2.84%  '

# Row 14
Set-CellText $ws 14 4 '15.33'
Set-CellText $ws 14 5 '  -1.56%  '

# Row 15
Set-CellText $ws 15 5 '  +1.68%  '

# Row 16
Set-CellText $ws 16 4 '2.580.50'
Set-CellText $ws 16 5 '  -1.72%  '

# Row 17
Set-CellText $ws 17 4 '2.255.03'
Set-CellText $ws 17 5 '  -1.32%  '

# Row 18
Set-CellText $ws 18 4 '43.000.90'
Set-CellText $ws 18 5 '  -0.55%  '

# Row 19
Set-CellText $ws 19 5 '  -1.19%  '

# Row 20
Set-CellText $ws 20 4 '6.70'
Set-CellText $ws 20 5 '  -2.30%  '

# Row 21
Set-CellText $ws 21 4 '71.59'
Set-CellText $ws 21 5 '  +0.00%  '

# Row 22
Set-CellText $ws 22 4 '2.36'
Set-CellText $ws 22 5 '  -6.47%  '

# Row 23
Set-CellText $ws 23 4 '231.09'
Set-CellText $ws 23 5 '  -0.72%  '

# Row 24
Set-CellText $ws 24 2 'PancakeSwap'
Set-CellText $ws 24 3 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
Set-CellText $ws 24 4 '2.89'
Set-CellText $ws 24 5 '  +0.24%  '

# Row 25
Set-CellText $ws 25 2 'InternetComputer(DFINITY)'
Set-CellText $ws 25 3 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-CellText $ws 25 4 '9.46'
Set-CellText $ws 25 5 '  -2.20%  '

# Row 26
Set-CellText $ws 26 5 '  +5.76%  '

# Row 27
Set-CellText $ws 27 5 '  +0.68%  '

# Row 28
Set-CellText $ws 28 4 '41.07'
Set-CellText $ws 28 5 '  +0.49%  '

# Row 29
Set-CellText $ws 29 2 'Toncoin'
Set-CellText $ws 29 3 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-CellText $ws 29 4 '2.24'
Set-CellText $ws 29 5 '  -0.44%  '

# Row 30
Set-CellText $ws 30 2 'WEMIXToken'
Set-CellText $ws 30 3 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
Set-CellText $ws 30 4 '3.29'
Set-CellText $ws 30 5 '  -1.40%  '

# Row 31
Set-CellText $ws 31 4 '172.44'
Set-CellText $ws 31 5 '  -0.07%  '

# Row 32
Set-CellText $ws 32 4 '21.17'
Set-CellText $ws 32 5 '  -1.25%  '

# Row 33
Set-CellText $ws 33 4 '0.0894'
Set-CellText $ws 33 5 '  -1.50%  '

# Row 34
Set-CellText $ws 34 4 '5.59'
Set-CellText $ws 34 5 '  -3.48%  '

# Row 35
Set-CellText $ws 35 4 '4.33'
Set-CellText $ws 35 5 '  +8.96%  '

# Row 36
Set-CellText $ws 36 5 '  +0.45%  '

# Row 37
Set-CellText $ws 37 2 'VeChain'
Set-CellText $ws 37 3 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-CellText $ws 37 4 '0.0371'
Set-CellText $ws 37 5 '  +3.62%  '

# Row 38
Set-CellText $ws 38 2 'RenderToken'
Set-CellText $ws 38 3 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-CellText $ws 38 4 '4.64'
Set-CellText $ws 38 5 '  -0.69%  '

# Row 39
Set-CellText $ws 39 5 '  +1.18%  '

# Row 40
Set-CellText $ws 40 5 '  -7.34%  '

# Row 41
Set-CellText $ws 41 2 'Celestia'
Set-CellText $ws 41 3 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-CellText $ws 41 4 '13.36'
Set-CellText $ws 41 5 '  -4.71%  '

# Row 42
Set-CellText $ws 42 2 'Algorand'
Set-CellText $ws 42 3 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-CellText $ws 42 4 '0.234'
Set-CellText $ws 42 5 '  -2.58%  '

# Row 43
Set-CellText $ws 43 4 '70.89'
Set-CellText $ws 43 5 '  -9.23%  '

# Row 44
Set-CellText $ws 44 5 '  +0.09%  '

# Row 45
Set-CellText $ws 45 5 '  -3.76%  '

# Row 46
Set-CellText $ws 46 4 '5.62'
Set-CellText $ws 46 5 '  -9.59%  '

# Row 47
Set-CellText $ws 47 4 '72.80'
Set-CellText $ws 47 5 '  +30.01%  '

# Row 48
Set-CellText $ws 48 2 'TheSandbox'
Set-CellText $ws 48 3 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-CellText $ws 48 4 '0.651'
Set-CellText $ws 48 5 '  +14.71%  '

# Row 49
Set-CellText $ws 49 2 'TrustWalletToken'
Set-CellText $ws 49 3 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-CellText $ws 49 4 '1.24'
Set-CellText $ws 49 5 '  -0.85%  '

# Row 50
Set-CellText $ws 50 2 'FraxShare'
Set-CellText $ws 50 3 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-CellText $ws 50 4 '8.39'
Set-CellText $ws 50 5 '  -3.84%  '

# Row 51
Set-CellText $ws 51 2 'Cronos'
Set-CellText $ws 51 3 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-CellText $ws 51 4 '0.0986'
Set-CellText $ws 51 5 '  -1.15%  '
